$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title / header text updates (shared string rich text runs) ---
$a8 = $ws.Range("A8")
$a8.Characters(21, 2).Text = "23"

$c9 = $ws.Range("C9")
$c9.Characters(47, 8).Text = "6/11/2023"
$c9.Characters(27, 9).Text = "6/5/2023"

# --- Numeric cell updates ---
$ws.Range("C16").Value = 4
$ws.Range("C17").Value = 10
$ws.Range("C18").Value = 2
$ws.Range("C19").Value = 11
$ws.Range("C20").Value = 12
$ws.Range("C21").Value = 39
$ws.Range("C24").Value = 25
$ws.Range("C25").Value = 13
$ws.Range("D16").Value = 4
$ws.Range("D17").Value = 9
$ws.Range("D18").Value = 3
$ws.Range("D19").Value = 12
$ws.Range("D20").Value = 5
$ws.Range("D21").Value = 33
$ws.Range("D24").Value = 29
$ws.Range("D25").Value = 6
$ws.Range("E16").Value = 0
$ws.Range("E17").Value = 11.111111111111
$ws.Range("E18").Value = -33.333333333333
$ws.Range("E19").Value = -8.333333333333
$ws.Range("E20").Value = 140
$ws.Range("E21").Value = 18.181818181818
$ws.Range("E24").Value = -13.793103448275
$ws.Range("E25").Value = 116.666666666667
$ws.Range("F16").Value = 14
$ws.Range("F17").Value = 36
$ws.Range("F18").Value = 8
$ws.Range("F19").Value = 50
$ws.Range("F20").Value = 23
$ws.Range("F21").Value = 131
$ws.Range("F24").Value = 112
$ws.Range("F25").Value = 48
$ws.Range("G16").Value = 12
$ws.Range("G17").Value = 26
$ws.Range("G18").Value = 13
$ws.Range("G19").Value = 50
$ws.Range("G20").Value = 12
$ws.Range("G21").Value = 117
$ws.Range("G24").Value = 112
$ws.Range("G25").Value = 41
$ws.Range("G27").Value = 4
$ws.Range("H15").Value = -100
$ws.Range("H16").Value = 16.666666666666
$ws.Range("H17").Value = 38.461538461538
$ws.Range("H18").Value = -38.461538461538
$ws.Range("H19").Value = 0
$ws.Range("H20").Value = 91.666666666666
$ws.Range("H21").Value = 11.965811965812
$ws.Range("H24").Value = 0
$ws.Range("H25").Value = 17.073170731707
$ws.Range("H26").Value = -100
$ws.Range("H27").Value = -25
$ws.Range("I16").Value = 86
$ws.Range("I17").Value = 147
$ws.Range("I18").Value = 44
$ws.Range("I19").Value = 272
$ws.Range("I20").Value = 101
$ws.Range("I21").Value = 660
$ws.Range("I24").Value = 568
$ws.Range("I25").Value = 247
$ws.Range("I27").Value = 25
$ws.Range("J16").Value = 85
$ws.Range("J17").Value = 128
$ws.Range("J18").Value = 62
$ws.Range("J19").Value = 314
$ws.Range("J20").Value = 112
$ws.Range("J21").Value = 713
$ws.Range("J24").Value = 623
$ws.Range("J25").Value = 199
$ws.Range("J28").Value = 7
$ws.Range("J29").Value = 7
$ws.Range("K16").Value = 1.176470588235
$ws.Range("K17").Value = 14.84375
$ws.Range("K18").Value = -29.032258064516
$ws.Range("K19").Value = -13.375796178343
$ws.Range("K20").Value = -9.821428571428
$ws.Range("K21").Value = -7.433380084151
$ws.Range("K24").Value = -8.828250401284
$ws.Range("K25").Value = 24.120603015075
$ws.Range("K27").Value = -3.846153846153
$ws.Range("K28").Value = -42.857142857142
$ws.Range("K29").Value = -71.428571428571
$ws.Range("L16").Value = 14.666666666666
$ws.Range("L17").Value = 27.826086956521
$ws.Range("L18").Value = 10
$ws.Range("L19").Value = 34.653465346534
$ws.Range("L20").Value = 57.8125
$ws.Range("L21").Value = 30.952380952381
$ws.Range("L24").Value = 56.473829201101
$ws.Range("L25").Value = 30.68783068783
$ws.Range("L27").Value = 25
$ws.Range("L28").Value = -33.333333333333
$ws.Range("L29").Value = -60
$ws.Range("M15").Value = -25
$ws.Range("M16").Value = -14.851485148514
$ws.Range("M17").Value = 126.153846153846
$ws.Range("M18").Value = -64.227642276422
$ws.Range("M19").Value = 85.034013605442
$ws.Range("M20").Value = -19.841269841269
$ws.Range("M21").Value = 14.782608695652
$ws.Range("M24").Value = 102.857142857143
$ws.Range("M25").Value = 12.272727272727
$ws.Range("N15").Value = -43.75
$ws.Range("N16").Value = -77.368421052631
$ws.Range("N17").Value = -7.54716981132
$ws.Range("N18").Value = -91.522157996146
$ws.Range("N19").Value = 8.8
$ws.Range("N20").Value = -93.239625167336
$ws.Range("N21").Value = -76.694915254237

# --- Cells changing to "0" shared-string display (t=s v=20, style 14) ---
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("F14").NumberFormat = "@"
$ws.Range("F14").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("F14").PasteSpecial(-4122)
$ws.Range("F15").NumberFormat = "@"
$ws.Range("F15").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("F15").PasteSpecial(-4122)
$ws.Range("F22").NumberFormat = "@"
$ws.Range("F22").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("F22").PasteSpecial(-4122)
$ws.Range("F26").NumberFormat = "@"
$ws.Range("F26").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("F26").PasteSpecial(-4122)

# --- Cells changing to "***.*" shared-string display (t=s v=21, style 14) ---
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("E26").PasteSpecial(-4122)
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("E27").PasteSpecial(-4122)
